$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)
$shape = $s.Shapes.Item(4)

# --- Update the "Write reusable scripts" bullet text, keeping it as a single
# run with its original formatting intact (operate on a precise character
# sub-range rather than TextRange.Paragraphs(), which avoids the run getting
# split into pieces when the replacement text partially overlaps the old one).
$tr = $shape.TextFrame.TextRange
$target = $tr.Paragraphs(3)
$sub = $tr.Characters($target.Start, $target.Length)
$sub.Text = "Write scripts that create maps"

# --- Insert a brand-new bullet paragraph right after it, inheriting the
# same paragraph/bullet + run formatting (InsertAfter with a leading carriage
# return starts a new paragraph rather than appending inline text).
$tr2 = $shape.TextFrame.TextRange
$para3 = $tr2.Paragraphs(3)
$para3.InsertAfter("`rUse Python to generate inputs")

# --- The text box auto-fits to its content (spAutoFit); pin the height to
# match the exact target extent from the authored edit. Shape.Height is a
# points value that gets truncated to EMU (1 pt = 12700 EMU) through a
# single-precision float, so nudge the literal fractionally above the exact
# quotient to land on the correct EMU after that truncation.
$shape.Height = 272.5761566322868
